# The commit swaps the two theme parts of the deck: the theme that was
# previously "Integral" (the one actually applied to the slide master /
# slides) becomes the stock "Office Theme" palette, while the palette that
# used to be the vanilla "Office Theme" (only ever linked from the notes
# master) becomes "Integral". Net visible effect for the deck's design is
# that the slide master/theme's 12 scheme colors change from the Integral
# greens/golds to the default Office blue/orange palette.
#
# PowerPoint's theme colors are exposed as a 12-slot ThemeColorScheme on the
# (Slide)Master's Theme; RGB values go in as COM "OLE_COLOR" longs, which
# are packed 0x00BBGGRR (blue in the high byte), not the usual 0xRRGGBB.

function ToOleColor([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" clrScheme, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (index 1-12).
$officeColors = @(
    @(0x00, 0x00, 0x00),   # 1  dk1
    @(0xFF, 0xFF, 0xFF),   # 2  lt1
    @(0x44, 0x54, 0x6A),   # 3  dk2
    @(0xE7, 0xE6, 0xE6),   # 4  lt2
    @(0x5B, 0x9B, 0xD5),   # 5  accent1
    @(0xED, 0x7D, 0x31),   # 6  accent2
    @(0xA5, 0xA5, 0xA5),   # 7  accent3
    @(0xFF, 0xC0, 0x00),   # 8  accent4
    @(0x44, 0x72, 0xC4),   # 9  accent5
    @(0x70, 0xAD, 0x47),   # 10 accent6
    @(0x05, 0x63, 0xC1),   # 11 hlink
    @(0x95, 0x4F, 0x72)    # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $scheme.Colors($i).RGB = ToOleColor $rgb[0] $rgb[1] $rgb[2]
}
